$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11 (row 1 header is unchanged).
# Columns A..AO (41 columns): League, Date, Time, Home, Away, then the odds columns.
$dataRows = @(
  @("Australian A-League Men", "2026-01-01", "01:00:00", "Auckland FC", "Newcastle Jets", 1.74, 1.76, 4.6, 4.8, 4.6, 4.7, 1.24, 1.03, 6.4, 1.16, 2.88, 1.48, 1.76, 2.2, 1.54, 2.68, 1.26, 2.3, 1000, 28, 1000, 110, 15, 12, 20, 1000, 15, 11, 17, 1000, 21, 17, 28, 1000, 7, 32),
  @("Australian A-League Men", "2026-01-01", "05:00:00", "Western Sydney Wanderers", "Macarthur FC", 1.93, 1.96, 4.2, 4.3, 3.95, 4.1, 1.31, 1.05, 4.7, 1.24, 2.3, 1.68, 1.52, 2.7, 1.66, 2.38, 1.3, 2.04, 23, 19.5, 36, 85, 14, 11, 18, 48, 14.5, 11.5, 17.5, 48, 23, 19.5, 30, 1000, 10.5, 38),
  @("Welsh Premiership", "2026-01-01", "09:30:00", "Colwyn Bay", "Flint Town United", 1.67, 2.08, 3.6, 6.2, 3.4, 8, 1.01, 1.04, 1.1, 1.23, 1.83, 1.43, 1.36, 2.46, 1.04, 1.04, 1.19, 1.92, 990, 990, 1000, 1000, 990, 990, 990, 1000, 1000, 990, 990, 1000, 1000, 1000, 1000, 1000, 1000, 1000),
  @("Saudi 1st Division", "2026-01-01", "09:35:00", "Al Orubah", "Al-Jndal", 1.02, 1000, 1.02, 1000, 1.02, 1000, 1.01, 1.01, 1.17, 1.01, 1.17, 1.01, 1.09, 1.39, 1.04, 1.04, 1.01, 1.01, 990, 990, 1000, 1000, 990, 990, 990, 1000, 1000, 990, 990, 1000, 1000, 1000, 1000, 1000, 1000, 1000),
  @("Saudi 1st Division", "2026-01-01", "12:00:00", "Al Faisaly ( KSA )", "Al-Raed (KSA)", 1.01, 1000, 1.01, 1000, 1.01, 950, 1.01, 1.01, 1.24, 1.01, 1.24, 1.01, 1.08, 1.02, 1.01, 1.01, 1.01, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000),
  @("English Premier League", "2026-01-01", "14:30:00", "Crystal Palace", "Fulham", 2.28, 2.3, 3.55, 3.65, 3.45, 3.5, 0, 1.08, 3.5, 1.37, 1.82, 2.18, 1.32, 3.95, 1.87, 2.06, 0, 0, 12.5, 13, 25, 70, 9.199999999999999, 7.4, 15, 46, 14, 11.5, 18.5, 55, 32, 26, 44, 130, 21, 50),
  @("English Premier League", "2026-01-01", "14:30:00", "Liverpool", "Leeds", 1.57, 1.58, 6.4, 6.6, 4.7, 4.8, 0, 1.04, 5.2, 1.23, 2.44, 1.68, 1.55, 2.76, 1.78, 2.24, 0, 0, 23, 27, 55, 200, 9.800000000000001, 10.5, 24, 70, 11, 9.4, 21, 70, 15, 15, 29, 90, 7, 80),
  @("Israeli Premier League", "2026-01-01", "15:30:00", "Beitar Jerusalem", "Hapoel Tel Aviv", 1.91, 2.12, 3.45, 4.5, 3.95, 4.7, 0, 0, 0, 0, 2.5, 1.52, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0),
  @("English Premier League", "2026-01-01", "17:00:00", "Brentford", "Tottenham", 2.28, 2.3, 3.35, 3.45, 3.65, 3.7, 0, 1.07, 3.9, 1.32, 1.96, 2.02, 1.38, 3.5, 1.79, 2.2, 0, 0, 15.5, 13.5, 23, 70, 10.5, 7.8, 14.5, 38, 15, 11.5, 18, 60, 32, 25, 40, 80, 18.5, 36),
  @("English Premier League", "2026-01-01", "17:00:00", "Sunderland", "Man City", 8.199999999999999, 8.4, 1.46, 1.47, 5, 5.1, 0, 1.05, 5, 1.23, 2.36, 1.7, 1.53, 2.8, 1.91, 2.02, 0, 0, 21, 9.4, 9, 13, 28, 11.5, 10, 14.5, 75, 30, 23, 36, 260, 120, 100, 120, 140, 6.4)
)

# Make sure the Date column (B) keeps its literal text value instead of being
# auto-converted into a date serial number by Excel.
$ws.Range("B2:B11").NumberFormat = "@"

for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $rowValues = $dataRows[$r]
    $excelRow = $r + 2   # data starts at row 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}
